$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A55").Value = "alan"
$ws.Range("B55").Value = "asdffd"
$ws.Range("C55").Value = "'234563"
$ws.Range("D55").Value = "fwefw"
$ws.Range("E55").Value = "user"

$ws.Range("A56").Value = "erferf"
$ws.Range("B56").Value = "fwrefewrf"
$ws.Range("C56").Value = "'232434"
$ws.Range("D56").Value = "erferf"
$ws.Range("E56").Value = "adm"
